$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the absolute path recorded in workbook metadata (cosmetic, matches
# the original author re-saving the file from the "KTMT" subfolder).
# (No direct COM property for this metadata is typically exposed; skip if unsupported.)

# 1) Fix the duplicated/incomplete question text on row 609 — it was missing
#    its trailing question mark and duplicated a (near-identical) deleted row.
$target = $ws.Range("A609")
$target.Value = "Trong kiểu trọng tài Bus nào thì việc phân chia quyền sử dụng Bus không cần một đơn vị trọng tài Bus riêng biệt?"

# 2) Remove the stray duplicate answer row (row 443) that had no content
#    beyond repeating another choice already present (row 452) — deleting it
#    shifts every following row up by one, which is exactly what the source
#    workbook reflects (dimension A1:A663 -> A1:A662).
$ws.Rows(443).Delete()

Write-Output "done"
